# feat: rebuild templates with native PPTX elements (editable text, shapes)
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# 1. Slide background -> solid dark fill (212121)
# ---------------------------------------------------------------------
$s.FollowMasterBackground = $false
$s.Background.Fill.Solid()
$s.Background.Fill.ForeColor.RGB = 0x212121

# ---------------------------------------------------------------------
# 2. Picture 1 -> Rounded Rectangle 1 (new translucent panel)
#    Duplicate "table_area" first so the new shape inherits the exact
#    same <p:style> block (lnRef/fillRef/effectRef/fontRef) that a
#    native PowerPoint auto-shape carries, then re-style/re-shape it.
# ---------------------------------------------------------------------
$tableArea = $s.Shapes.Item(3)
$rr = $tableArea.Duplicate().Item(1)
$rr.ZOrder(1)  # msoSendToBack -> becomes the first shape in the tree

$rr.Name = "Rounded Rectangle 1"
$rr.AutoShapeType = 5  # msoShapeRoundedRectangle

$rr.Left = -82570 / 12700
$rr.Top = -82570 / 12700
$rr.Width = 24549171 / 12700
$rr.Height = 13881232 / 12700

$rr.Fill.Visible = $true
$rr.Fill.ForeColor.RGB = 0x212121
$rr.Fill.Transparency = 0.95

$rr.Line.Visible = $true
$rr.Line.ForeColor.RGB = 0xDEE0D6
$rr.Line.Weight = 4

# Remove the now-redundant picture
$s.Shapes.Item(2).Delete()

# ---------------------------------------------------------------------
# 3. slide_title -> TextBox 2  (big centred glyph placeholder)
# ---------------------------------------------------------------------
$tb2 = $s.Shapes.Item(2)
$tb2.Name = "TextBox 2"

$tr2 = $tb2.TextFrame.TextRange
$tr2.Font.Size = 77.3
$tr2.Font.Bold = $false
$tr2.Font.Italic = $false
$tr2.Font.Name = "Rajdhani"
$tr2.Font.Color.RGB = 0x000000
$tr2.Text = "_x0016_"
$tr2.ParagraphFormat.Alignment = 2  # ppAlignCenter

$tb2.Left = 3445824 / 12700
$tb2.Top = 2616189 / 12700
$tb2.Width = 17492380 / 12700
$tb2.Height = 9699223 / 12700

# ---------------------------------------------------------------------
# 4. table_area -> TextBox 3  (source footnote)
# ---------------------------------------------------------------------
$tb3 = $s.Shapes.Item(3)
$tb3.Name = "TextBox 3"
$tb3.Line.Visible = $false

$tr3 = $tb3.TextFrame.TextRange
$tr3.Font.Size = 20
$tr3.Font.Bold = $false
$tr3.Font.Italic = $false
$tr3.Font.Name = "Quicksand (TT)"
$tr3.Font.Color.RGB = 0x000000
$tr3.Text = "Source: FastTrack as of 12/31/2025. Dataset inception: 9/1/1988. "
$tr3.ParagraphFormat.Alignment = 2  # ppAlignCenter

$tb3.TextFrame.WordWrap = $true
$tb3.TextFrame.AutoSize = 1  # ppAutoSizeShapeToFitText -> <a:spAutoFit/>

$tb3.Left = 6458955 / 12700
$tb3.Top = 12827020 / 12700
$tb3.Width = 11466118 / 12700
$tb3.Height = 446135 / 12700

# ---------------------------------------------------------------------
# 5. footnote -> TextBox 4  (headline)
# ---------------------------------------------------------------------
$tb4 = $s.Shapes.Item(4)
$tb4.Name = "TextBox 4"

$tr4 = $tb4.TextFrame.TextRange
$tr4.Font.Size = 36
$tr4.Font.Bold = $false
$tr4.Font.Italic = $false
$tr4.Font.Name = "Rajdhani"
$tr4.Font.Color.RGB = 0x212121
$tr4.Text = "IF YOU WANT THE S&P 500, JUST BUY THE S&P 500!"
$tr4.ParagraphFormat.Alignment = 1  # ppAlignLeft

$tb4.Left = 1066830 / 12700
$tb4.Top = 1143000 / 12700
$tb4.Width = 9829800 / 12700
$tb4.Height = 559155 / 12700
